$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1) changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2) changed
$ws.Range("B2").Value = -10.848549352678246
$ws.Range("C2").Value = 10.722787820974062
$ws.Range("D2").Value = 3.3085748382100282
$ws.Range("E2").Value = 3.027274908548629

# Row 3 data values (B3:E3) changed
$ws.Range("B3").Value = 9.1606185307708188
$ws.Range("C3").Value = 13.270904995222017
$ws.Range("D3").Value = 18.036638296737841
$ws.Range("E3").Value = -2.1551504644406752

# Update the active selection to match the edited range
$ws.Range("B1:E3").Select()
